$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 209, shifting existing rows 209-253 down to 210-254
$ws.Rows.Item(209).Insert()

# Populate the new row 209 with the new data
$ws.Cells.Item(209, 1).Value = 4
$ws.Cells.Item(209, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(209, 3).Value = "Los Lagos"
$ws.Cells.Item(209, 4).Value = 44785
$ws.Cells.Item(209, 5).Value = 10
$ws.Cells.Item(209, 6).Value = "Fruta"
$ws.Cells.Item(209, 7).Value = 100109
$ws.Cells.Item(209, 8).Value = "Uva"
$ws.Cells.Item(209, 9).Value = 100109001
$ws.Cells.Item(209, 10).Value = "Uva"
$ws.Cells.Item(209, 11).Value = "Red Globe"
$ws.Cells.Item(209, 12).Value = "Primera"
$ws.Cells.Item(209, 13).Value = 300
$ws.Cells.Item(209, 14).Value = 15000
$ws.Cells.Item(209, 15).Value = 16000
$ws.Cells.Item(209, 16).Value = 15500
$ws.Cells.Item(209, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(209, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(209, 19).Value = 1938
$ws.Cells.Item(209, 20).Value = 8
